# Bug fixes on test 5 generator: update Write Latency "min" (col O) and
# "max" (col P) figures for rows 3-15 on the active sheet.
#
# The source cells store their values as text (not numbers), even when the
# text looks numeric (e.g. "644"). Assigning a numeric-looking string
# straight to Range.Value would make Excel auto-coerce it into a real
# number, which would change the cell's stored type. To avoid that we
# stage each value in a scratch cell that has been explicitly formatted as
# Text ("@"), copy it, and paste-special *values only* into the real
# target cell. That keeps the destination's existing style/number format
# untouched while still writing a genuine text value, exactly like the
# values already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell well outside the populated A1:X15 table; cleared after each
# use so it never lingers in the sheet's used range.
$helper = $ws.Range("Z1")

$updates = @(
  @{Cell="O3";  Value="1839"},
  @{Cell="P3";  Value="12516k"},
  @{Cell="O4";  Value="1806"},
  @{Cell="P4";  Value="6480.3k"},
  @{Cell="O5";  Value="1361"},
  @{Cell="P5";  Value="8460.8k"},
  @{Cell="O6";  Value="2"},
  @{Cell="P6";  Value="6680"},
  @{Cell="O7";  Value="1775"},
  @{Cell="P7";  Value="1696.7k"},
  @{Cell="O8";  Value="1615"},
  @{Cell="P8";  Value="1079.8k"},
  @{Cell="O9";  Value="2"},
  @{Cell="P9";  Value="2396"},
  @{Cell="O10"; Value="2"},
  @{Cell="P10"; Value="6141"},
  @{Cell="O11"; Value="1907"},
  @{Cell="P11"; Value="1807.9k"},
  @{Cell="O12"; Value="3"},
  @{Cell="P12"; Value="1280"},
  @{Cell="O13"; Value="1860"},
  @{Cell="P13"; Value="6696.6k"},
  @{Cell="O14"; Value="1888"},
  @{Cell="P14"; Value="2682.8k"},
  @{Cell="O15"; Value="1337"},
  @{Cell="P15"; Value="3111.3k"}
)

foreach ($update in $updates) {
  $helper.NumberFormat = "@"
  $helper.Value = $update.Value
  $helper.Copy()
  $ws.Range($update.Cell).PasteSpecial(-4163)
  $helper.Clear()
}
